# Updates the cryptos list sheet with refreshed price/volume data,
# a new OKB row inserted at row 8 (pushing later rows down by one),
# and the final row (Cronos) falling off the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin, Link, Price, "Volume(1h)"
$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.787.33', '  +2.15%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.860.02', '  +1.69%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9995', '  +0.05%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '245.16', '  +0.96%  '),
    @(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6413', '  +3.46%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9999', '  -0.04%  '),
    @(8, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '47.58', '  +4.61%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07547', '  +2.77%  '),
    @(10, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2979', '  +2.57%  '),
    @(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.54', '  +5.55%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07681', '  +0.67%  '),
    @(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.860.67', '  +1.58%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.045', '  +1.71%  '),
    @(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6923', '  +3.51%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '84.04', '  +2.00%  '),
    @(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009855', '  +9.77%  '),
    @(18, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.112', '  +4.77%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.797.63', '  +2.21%  '),
    @(20, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.114.28', '  +1.40%  '),
    @(21, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '236.17', '  +0.15%  '),
    @(22, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.67', '  +1.81%  '),
    @(23, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.00%  '),
    @(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.515', '  +2.34%  '),
    @(25, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  -0.03%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '159.21', '  +0.45%  '),
    @(27, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1421', '  +2.21%  '),
    @(28, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.559', '  +0.49%  '),
    @(29, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.95', '  +1.93%  '),
    @(30, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06193', '  +6.16%  '),
    @(31, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.494', '  +0.57%  '),
    @(32, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.293', '  +5.94%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.162', '  +1.93%  '),
    @(34, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.108', '  +0.79%  '),
    @(35, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.899', '  +2.41%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.176', '  +3.39%  '),
    @(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7315', '  +0.97%  '),
    @(38, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.603', '  -0.12%  '),
    @(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.814', '  -1.23%  '),
    @(40, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01786', '  +1.53%  '),
    @(41, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.209.84', '  -1.44%  '),
    @(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.291', '  +1.12%  '),
    @(43, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9204', '  +1.41%  '),
    @(44, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.000', '  +0.03%  '),
    @(45, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.021.61', '  +1.66%  '),
    @(46, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '102.10', '  +0.30%  '),
    @(47, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '66.64', '  +1.34%  '),
    @(48, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000118', '  +1.39%  '),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.218', '  +1.16%  '),
    @(50, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4063', '  +0.75%  '),
    @(51, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.676', '  +5.99%  ')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[1]
    $ws.Range("C$rowNum").Value = $r[2]
    # Leading apostrophe forces text so values like "1.000" or "0.9999"
    # are not auto-converted to numbers by Excel.
    $ws.Range("D$rowNum").Value = "'" + $r[3]
    $ws.Range("E$rowNum").Value = $r[4]
}